$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 1061
$ws1.Range("F16").Value = 93
$ws1.Range("F24").Value = 399
$ws1.Range("F25").Value = 171
$ws1.Range("F32").Value = 449
$ws1.Range("F33").Value = 449
$ws1.Range("F41").Value = 1252
$ws1.Range("F42").Value = 3254
$ws1.Range("F49").Value = 473

# Sheet "本地生活" (local life) - update F column values
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1999
$ws3.Range("F7").Value = 557
$ws3.Range("F10").Value = 1059

# Sheet "全部类型" (all types) - update F column values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1999
$ws4.Range("F6").Value = 557
$ws4.Range("F8").Value = 1059
$ws4.Range("F16").Value = 93
$ws4.Range("F24").Value = 399
$ws4.Range("F25").Value = 171
$ws4.Range("F31").Value = 449
$ws4.Range("F32").Value = 450
$ws4.Range("F40").Value = 3254
